$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2035.0476
$ws.Range("I15").Value = 2035.0476
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 6105.142800000001
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -5936.142800000001

$ws.Range("H33").Value = 728.45
$ws.Range("I33").Value = 520.625
$ws.Range("J33").Value = 1559.75
$ws.Range("K33").Value = 520.625
$ws.Range("L33").Value = 1559.75
$ws.Range("M33").Value = -291.625
$ws.Range("N33").Value = -2017.75

$ws.Range("H98").Value = 20321.787
$ws.Range("I98").Value = 21562.037
$ws.Range("J98").Value = 14740.667
$ws.Range("K98").Value = 21562.037
$ws.Range("L98").Value = 14740.667
$ws.Range("M98").Value = -20064.037
$ws.Range("N98").Value = -17736.667

$ws.Range("H116").Value = 1513895.6
$ws.Range("I116").Value = 1593048.2
$ws.Range("J116").Value = 9996
$ws.Range("K116").Value = 1593048.2
$ws.Range("L116").Value = 9996
$ws.Range("M116").Value = -1589606.2
$ws.Range("N116").Value = -16880

$ws.Range("H122").Value = 20321.787
$ws.Range("I122").Value = 21562.037
$ws.Range("J122").Value = 14740.667
$ws.Range("K122").Value = 64686.111
$ws.Range("L122").Value = 44222.001
$ws.Range("M122").Value = -62236.111
$ws.Range("N122").Value = -49122.001

$ws.Range("H137").Value = 37316
$ws.Range("I137").Value = 43779.2
$ws.Range("J137").Value = 5000
$ws.Range("K137").Value = 131337.6
$ws.Range("L137").Value = 15000
$ws.Range("M137").Value = -128787.6
$ws.Range("N137").Value = -20100

$ws.Range("H141").Value = 3726.3235
$ws.Range("I141").Value = 3501.8276
$ws.Range("J141").Value = 5028.4
$ws.Range("K141").Value = 10505.4828
$ws.Range("L141").Value = 15085.2
$ws.Range("M141").Value = -5325.4828
$ws.Range("N141").Value = -25445.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 2718.7144
$ws.Range("I5").Value = 2718.7144
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 2718.7144
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -2606.7144
$ws.Range("N5").ClearContents()

$ws.Range("H32").Value = 10921.108
$ws.Range("I32").Value = 11714.576
$ws.Range("J32").Value = 4375
$ws.Range("K32").Value = 11714.576
$ws.Range("L32").Value = 4375
$ws.Range("M32").Value = -11427.576
$ws.Range("N32").Value = -4949

$ws.Range("H61").Value = 4432.978
$ws.Range("I61").Value = 4641.675
$ws.Range("J61").Value = 3041.6667
$ws.Range("K61").Value = 4641.675
$ws.Range("L61").Value = 3041.6667
$ws.Range("M61").Value = -4429.675
$ws.Range("N61").Value = -3465.6667

$ws.Range("H63").Value = 2010.5
$ws.Range("I63").Value = 2034.8572
$ws.Range("J63").Value = 1840
$ws.Range("K63").Value = 2034.8572
$ws.Range("L63").Value = 1840
$ws.Range("M63").Value = -1348.8572
$ws.Range("N63").Value = -3212

$ws.Range("H66").Value = 2010.5
$ws.Range("I66").Value = 2034.8572
$ws.Range("J66").Value = 1840
$ws.Range("K66").Value = 10174.286
$ws.Range("L66").Value = 9200
$ws.Range("M66").Value = -6742.286
$ws.Range("N66").Value = -16064

$ws.Range("H104").Value = 48998.5
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 48998.5
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 48998.5
$ws.Range("N104").Value = -55986.5

$ws.Range("H132").Value = 6356.5854
$ws.Range("I132").Value = 7262.231
$ws.Range("J132").Value = 4786.8
$ws.Range("K132").Value = 21786.693
$ws.Range("L132").Value = 14360.4
$ws.Range("M132").Value = -19256.693
$ws.Range("N132").Value = -19420.4

$ws.Range("H136").Value = 4432.978
$ws.Range("I136").Value = 4641.675
$ws.Range("J136").Value = 3041.6667
$ws.Range("K136").Value = 13925.025
$ws.Range("L136").Value = 9125.000100000001
$ws.Range("M136").Value = -11375.025
$ws.Range("N136").Value = -14225.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 2718.7144
$ws.Range("I4").Value = 2718.7144
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 2718.7144
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -2603.7144
$ws.Range("N4").ClearContents()

$ws.Range("H22").Value = 631
$ws.Range("I22").Value = 631
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 631
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -458
$ws.Range("N22").ClearContents()

$ws.Range("H99").Value = 42401.5
$ws.Range("I99").Value = 72823
$ws.Range("J99").Value = 11980
$ws.Range("K99").Value = 72823
$ws.Range("L99").Value = 11980
$ws.Range("M99").Value = -71325
$ws.Range("N99").Value = -14976

$ws.Range("H123").Value = 72555
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 72555
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 72555
$ws.Range("N123").Value = -82355

$ws.Range("H134").Value = 14124.7
$ws.Range("I134").Value = 16158.529
$ws.Range("J134").Value = 2599.6667
$ws.Range("K134").Value = 48475.587
$ws.Range("L134").Value = 7799.000100000001
$ws.Range("M134").Value = -45940.587
$ws.Range("N134").Value = -12869.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 6428.2812
$ws.Range("I7").Value = 8433.708000000001
$ws.Range("J7").Value = 412
$ws.Range("K7").Value = 8433.708000000001
$ws.Range("L7").Value = 412
$ws.Range("M7").Value = -8320.708000000001
$ws.Range("N7").Value = -638

$ws.Range("H31").Value = 2965.3967
$ws.Range("I31").Value = 2956.9546
$ws.Range("J31").Value = 2984.9473
$ws.Range("K31").Value = 2956.9546
$ws.Range("L31").Value = 2984.9473
$ws.Range("M31").Value = -2661.9546
$ws.Range("N31").Value = -3574.9473

$ws.Range("H34").Value = 2965.3967
$ws.Range("I34").Value = 2956.9546
$ws.Range("J34").Value = 2984.9473
$ws.Range("K34").Value = 2956.9546
$ws.Range("L34").Value = 2984.9473
$ws.Range("M34").Value = -2754.9546
$ws.Range("N34").Value = -3388.9473

$ws.Range("H54").Value = 10139.5
$ws.Range("I54").Value = 279
$ws.Range("J54").Value = 20000
$ws.Range("K54").Value = 279
$ws.Range("L54").Value = 20000
$ws.Range("M54").Value = 379
$ws.Range("N54").Value = -21316

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 28251670
$ws.Range("I4").Value = 51219604
$ws.Range("J4").Value = 13178961
$ws.Range("K4").Value = 153658812
$ws.Range("L4").Value = 39536883
$ws.Range("M4").Value = -153658700
$ws.Range("N4").Value = -39537107

$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()

$ws.Range("H17").Value = 1699.7333
$ws.Range("I17").Value = 1211
$ws.Range("J17").Value = 1774.9231
$ws.Range("K17").Value = 3633
$ws.Range("L17").Value = 5324.7693
$ws.Range("M17").Value = -3464
$ws.Range("N17").Value = -5662.7693

$ws.Range("H34").Value = 587927.8
$ws.Range("I34").Value = 1590216.1
$ws.Range("J34").Value = 3259.5833
$ws.Range("K34").Value = 4770648.300000001
$ws.Range("L34").Value = 9778.749899999999
$ws.Range("M34").Value = -4770564.300000001
$ws.Range("N34").Value = -9946.749899999999

$ws.Range("H39").Value = 1470
$ws.Range("I39").Value = 300
$ws.Range("J39").Value = 12000
$ws.Range("K39").Value = 900
$ws.Range("L39").Value = 36000
$ws.Range("M39").Value = -606
$ws.Range("N39").Value = -36588

$ws.Range("H55").Value = 4700.7144
$ws.Range("I55").Value = 1281.6
$ws.Range("J55").Value = 5769.1875
$ws.Range("K55").Value = 3844.8
$ws.Range("L55").Value = 17307.5625
$ws.Range("M55").Value = -3667.8
$ws.Range("N55").Value = -17661.5625

$ws.Range("H98").Value = 943.5714
$ws.Range("I98").Value = 994.7
$ws.Range("J98").Value = 815.75
$ws.Range("K98").Value = 2984.1
$ws.Range("L98").Value = 2447.25
$ws.Range("M98").Value = -1486.1
$ws.Range("N98").Value = -5443.25

$ws.Range("H122").Value = 2038.6666
$ws.Range("I122").Value = 488.55554
$ws.Range("J122").Value = 2555.3704
$ws.Range("K122").Value = 4396.99986
$ws.Range("L122").Value = 22998.3336
$ws.Range("M122").Value = -1946.99986
$ws.Range("N122").Value = -27898.3336

$ws.Range("H129").Value = 5557279
$ws.Range("I129").Value = 1087
$ws.Range("J129").Value = 10002233
$ws.Range("K129").Value = 3261
$ws.Range("L129").Value = 30006699
$ws.Range("M129").Value = 1739
$ws.Range("N129").Value = -30016699

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 22535.715
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 22535.715
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 22535.715
$ws.Range("N92").Value = -26279.715

$ws.Range("H102").Value = 6279.718
$ws.Range("I102").Value = 6614.1763
$ws.Range("J102").Value = 4005.4
$ws.Range("K102").Value = 6614.1763
$ws.Range("L102").Value = 4005.4
$ws.Range("M102").Value = -4992.1763
$ws.Range("N102").Value = -7249.4

$ws.Range("H107").Value = 331.1875
$ws.Range("I107").Value = 413.33334
$ws.Range("J107").Value = 84.75
$ws.Range("K107").Value = 413.33334
$ws.Range("L107").Value = 84.75
$ws.Range("M107").Value = 1506.66666
$ws.Range("N107").Value = -3924.75

$ws.Range("H117").Value = 36543
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 36543
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 36543
$ws.Range("N117").Value = -43427

$ws.Range("H132").Value = 3947.3877
$ws.Range("I132").Value = 3461.2563
$ws.Range("J132").Value = 5843.3
$ws.Range("K132").Value = 10383.7689
$ws.Range("L132").Value = 17529.9
$ws.Range("M132").Value = -7853.768899999999
$ws.Range("N132").Value = -22589.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4389.5
$ws.Range("I46").Value = 920
$ws.Range("J46").Value = 6124.25
$ws.Range("K46").Value = 920
$ws.Range("L46").Value = 6124.25
$ws.Range("M46").Value = -732
$ws.Range("N46").Value = -6500.25

$ws.Range("H122").Value = 9470.75
$ws.Range("I122").Value = 9968.091
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 29904.273
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -27454.273
$ws.Range("N122").Value = -16900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 12987.526
$ws.Range("I132").Value = 14301.8545
$ws.Range("J132").Value = 5977.778
$ws.Range("K132").Value = 42905.5635
$ws.Range("L132").Value = 17933.334
$ws.Range("M132").Value = -40375.5635
$ws.Range("N132").Value = -22993.334

$ws.Range("H136").Value = 227416.83
$ws.Range("I136").Value = 270721.75
$ws.Range("J136").Value = 3018.6365
$ws.Range("K136").Value = 812165.25
$ws.Range("L136").Value = 9055.9095
$ws.Range("M136").Value = -809615.25
$ws.Range("N136").Value = -14155.9095
